# feat: removes processor_full_name (#100)
#
# The first column of the sheet ("processor_full_name", containing
# the values "Anna Apple" / "John Smith" / "Anna Apple") is removed
# entirely. All the other columns shift one position to the left
# (what used to be column B is now column A, ... column K is now J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "processor_full_name" column; Excel automatically
# shifts every remaining column one place to the left.
$ws.Columns("A:A").Delete()

# Re-apply a plain "Normal" style to every surviving cell that holds
# data (mirrors how the source workbook marks its populated cells).
$used = $ws.UsedRange
foreach ($area in $used.SpecialCells(2).Areas) {
    foreach ($cell in $area.Cells) {
        $cell.Style = "Normal"
    }
}

# The active selection moves back to the top-left cell of the sheet.
[void]$ws.Range("A1").Select()
